$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "244.36"
Set-TextValue $ws.Range("D3") "23.88"
Set-TextValue $ws.Range("D4") "5.300"
Set-TextValue $ws.Range("D5") "0.05880"
Set-TextValue $ws.Range("D6") "6.479"
Set-TextValue $ws.Range("D7") "3.331"
Set-TextValue $ws.Range("D8") "0.8176"
Set-TextValue $ws.Range("D9") "0.8928"
Set-TextValue $ws.Range("D10") "0.1389"
Set-TextValue $ws.Range("D11") "0.07238"
Set-TextValue $ws.Range("D12") "0.03100"
Set-TextValue $ws.Range("D14") "0.09355"
Set-TextValue $ws.Range("D15") "3.820"
Set-TextValue $ws.Range("D17") "0.04730"
Set-TextValue $ws.Range("D18") "0.0006020"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue $ws.Range("D19") "0.006216"
$ws.Range("E19").Value = "18TigerCashTCH"
Set-TextValue $ws.Range("D20") "0.001261"
Set-TextValue $ws.Range("D21") "0.004618"
Set-TextValue $ws.Range("D22") "0.00008700"
Set-TextValue $ws.Range("D23") "3.553"
Set-TextValue $ws.Range("D24") "2.177"
Set-TextValue $ws.Range("D25") "0.3200"
Set-TextValue $ws.Range("D40") "0.03807"
Set-TextValue $ws.Range("D41") "0.006320"
Set-TextValue $ws.Range("D42") "0.1058"
Set-TextValue $ws.Range("D43") "0.002532"
$ws.Range("E43").Value = "42CEJICEJI"
Set-TextValue $ws.Range("D44") "0.007100"
Set-TextValue $ws.Range("D45") "0.00005366"
Set-TextValue $ws.Range("D47") "0.5399"
Set-TextValue $ws.Range("D48") "0.01835"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
Set-TextValue $ws.Range("D49") "0.00002100"
Set-TextValue $ws.Range("D50") "0.0002000"
